$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap home/away match data (columns F:V) between row pairs that share the same kickoff time (column E) ---
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
$pairs = @(
    @(39, 40),
    @(51, 52),
    @(59, 60),
    @(64, 65),
    @(82, 83),
    @(85, 87),
    @(97, 98),
    @(99, 100),
    @(124, 125),
    @(141, 142),
    @(169, 170),
    @(173, 174),
    @(184, 185),
    @(196, 197),
    @(207, 208),
    @(214, 215)
)

foreach ($pair in $pairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]
    foreach ($c in $cols) {
        $addrA = "${c}${rowA}"
        $addrB = "${c}${rowB}"
        $valA = $ws.Range($addrA).Value()
        $valB = $ws.Range($addrB).Value()
        $ws.Range($addrA).Value = $valB
        $ws.Range($addrB).Value = $valA
    }
}

# --- Append 4 new match rows (222:225), carrying the A/E column styles down from row 221 ---
$ws.Range("A221").Copy($ws.Range("A222:A225"))
$ws.Range("E221").Copy($ws.Range("E222:E225"))

# Row 222
$ws.Cells.Item(222, 1).Value = 221
$ws.Cells.Item(222, 2).Value = "spain"
$ws.Cells.Item(222, 3).Value = "laliga2"
$ws.Cells.Item(222, 4).Value = "2023-2024"
$ws.Cells.Item(222, 5).Value = 45279.79166666666
$ws.Cells.Item(222, 6).Value = "Albacete"
$ws.Cells.Item(222, 7).Value = 1
$ws.Cells.Item(222, 8).Value = "Eldense"
$ws.Cells.Item(222, 9).Value = 1
$ws.Cells.Item(222, 10).Value = 1.61
$ws.Cells.Item(222, 11).Value = "16/12/2023 21:12"
$ws.Cells.Item(222, 12).Value = 1.66
$ws.Cells.Item(222, 13).Value = "19/12/2023 18:59"
$ws.Cells.Item(222, 14).Value = 3.91
$ws.Cells.Item(222, 15).Value = "16/12/2023 21:12"
$ws.Cells.Item(222, 16).Value = 3.8
$ws.Cells.Item(222, 17).Value = "19/12/2023 18:59"
$ws.Cells.Item(222, 18).Value = 5.98
$ws.Cells.Item(222, 19).Value = "16/12/2023 21:12"
$ws.Cells.Item(222, 20).Value = 6.04
$ws.Cells.Item(222, 21).Value = "19/12/2023 18:59"
$ws.Cells.Item(222, 22).Value = "https://www.betexplorer.com/football/spain/laliga2/albacete-eldense/AZNvRg2C/"

# Row 223
$ws.Cells.Item(223, 1).Value = 222
$ws.Cells.Item(223, 2).Value = "spain"
$ws.Cells.Item(223, 3).Value = "laliga2"
$ws.Cells.Item(223, 4).Value = "2023-2024"
$ws.Cells.Item(223, 5).Value = 45279.79166666666
$ws.Cells.Item(223, 6).Value = "Racing Santander"
$ws.Cells.Item(223, 7).Value = 2
$ws.Cells.Item(223, 8).Value = "Andorra"
$ws.Cells.Item(223, 9).Value = 0
$ws.Cells.Item(223, 10).Value = 2.42
$ws.Cells.Item(223, 11).Value = "16/12/2023 18:13"
$ws.Cells.Item(223, 12).Value = 2.07
$ws.Cells.Item(223, 13).Value = "19/12/2023 18:56"
$ws.Cells.Item(223, 14).Value = 3.23
$ws.Cells.Item(223, 15).Value = "16/12/2023 18:13"
$ws.Cells.Item(223, 16).Value = 3.41
$ws.Cells.Item(223, 17).Value = "19/12/2023 18:56"
$ws.Cells.Item(223, 18).Value = 3.16
$ws.Cells.Item(223, 19).Value = "16/12/2023 18:13"
$ws.Cells.Item(223, 20).Value = 3.97
$ws.Cells.Item(223, 21).Value = "19/12/2023 18:56"
$ws.Cells.Item(223, 22).Value = "https://www.betexplorer.com/football/spain/laliga2/racing-santander-fc-andorra/OpS8mYe0/"

# Row 224
$ws.Cells.Item(224, 1).Value = 223
$ws.Cells.Item(224, 2).Value = "spain"
$ws.Cells.Item(224, 3).Value = "laliga2"
$ws.Cells.Item(224, 4).Value = "2023-2024"
$ws.Cells.Item(224, 5).Value = 45279.89583333334
$ws.Cells.Item(224, 6).Value = "Elche"
$ws.Cells.Item(224, 7).Value = 0
$ws.Cells.Item(224, 8).Value = "Mirandes"
$ws.Cells.Item(224, 9).Value = 0
$ws.Cells.Item(224, 10).Value = 1.74
$ws.Cells.Item(224, 11).Value = "16/12/2023 21:12"
$ws.Cells.Item(224, 12).Value = 1.92
$ws.Cells.Item(224, 13).Value = "19/12/2023 21:13"
$ws.Cells.Item(224, 14).Value = 3.76
$ws.Cells.Item(224, 15).Value = "16/12/2023 21:12"
$ws.Cells.Item(224, 16).Value = 3.51
$ws.Cells.Item(224, 17).Value = "19/12/2023 21:13"
$ws.Cells.Item(224, 18).Value = 4.98
$ws.Cells.Item(224, 19).Value = "16/12/2023 21:12"
$ws.Cells.Item(224, 20).Value = 4.49
$ws.Cells.Item(224, 21).Value = "19/12/2023 21:13"
$ws.Cells.Item(224, 22).Value = "https://www.betexplorer.com/football/spain/laliga2/elche-mirandes/IaRnPXWO/"

# Row 225
$ws.Cells.Item(225, 1).Value = 224
$ws.Cells.Item(225, 2).Value = "spain"
$ws.Cells.Item(225, 3).Value = "laliga2"
$ws.Cells.Item(225, 4).Value = "2023-2024"
$ws.Cells.Item(225, 5).Value = 45279.89583333334
$ws.Cells.Item(225, 6).Value = "Espanyol"
$ws.Cells.Item(225, 7).Value = 3
$ws.Cells.Item(225, 8).Value = "Burgos CF"
$ws.Cells.Item(225, 9).Value = 3
$ws.Cells.Item(225, 10).Value = 1.73
$ws.Cells.Item(225, 11).Value = "16/12/2023 18:42"
$ws.Cells.Item(225, 12).Value = 1.68
$ws.Cells.Item(225, 13).Value = "19/12/2023 21:20"
$ws.Cells.Item(225, 14).Value = 3.69
$ws.Cells.Item(225, 15).Value = "16/12/2023 18:42"
$ws.Cells.Item(225, 16).Value = 3.75
$ws.Cells.Item(225, 17).Value = "19/12/2023 21:27"
$ws.Cells.Item(225, 18).Value = 5.15
$ws.Cells.Item(225, 19).Value = "16/12/2023 18:42"
$ws.Cells.Item(225, 20).Value = 6.01
$ws.Cells.Item(225, 21).Value = "19/12/2023 21:27"
$ws.Cells.Item(225, 22).Value = "https://www.betexplorer.com/football/spain/laliga2/espanyol-burgos-cf/lKe1reXn/"

